# The workbook's single sheet contains a daily-price dataset ("Poroto granado")
# with one header row (row 1) and data rows 2-200. The edit inserts one new
# data record as a new row 127, pushing the former rows 127-200 down to 128-201.
#
# All rows share the same constant values in columns A, B, C, E, F, G, H, I, R
# (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría, Variedad,
# Calidad, Clasificación), so the new row reuses those constants and only
# supplies fresh values for D (Fecha), J (Volumen), K (Precio mínimo),
# L (Precio máximo), M (Precio promedio ponderado), N (Unidad de
# comercialización), O (Origen), P (Precio $/Kg) and Q (Kg o Unidades).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127; this shifts rows 127:200 down to 128:201
# and expands the used range / dimension to A1:R201 automatically.
$ws.Rows.Item(127).Insert()

$newRow = 127

# Columns that are constant across every data row in this dataset.
$ws.Cells.Item($newRow, 1).Value2 = 9
$ws.Cells.Item($newRow, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($newRow, 3).Value2 = "Metropolitana"
$ws.Cells.Item($newRow, 5).Value2 = 13
$ws.Cells.Item($newRow, 6).Value2 = 100112030
$ws.Cells.Item($newRow, 7).Value2 = "Poroto granado"
$ws.Cells.Item($newRow, 8).Value2 = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value2 = "Primera"
$ws.Cells.Item($newRow, 18).Value2 = "Hortaliza"

# New record's own data.
$ws.Cells.Item($newRow, 4).Value2 = 44572
$ws.Cells.Item($newRow, 10).Value2 = 79
$ws.Cells.Item($newRow, 11).Value2 = 25000
$ws.Cells.Item($newRow, 12).Value2 = 26000
$ws.Cells.Item($newRow, 13).Value2 = 25506
$ws.Cells.Item($newRow, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item($newRow, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item($newRow, 16).Value2 = 1020
$ws.Cells.Item($newRow, 17).Value2 = 25

# Column D uses a date-time number format across the sheet (style index 2);
# make sure the newly inserted cell keeps that format (Insert() already
# copies the formatting from the row above, but set it explicitly to be safe).
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat

Write-Host "Done. New dimension: $($ws.UsedRange.Address())"
